$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = 12
$ws.Cells.Item(2, 2).Value = "h`$_{q}`$"
$ws.Cells.Item(2, 3).Value = 0.001750130549641798
$ws.Cells.Item(3, 1).Value = 37
$ws.Cells.Item(3, 2).Value = "`$F_{q}`$"
$ws.Cells.Item(3, 3).Value = 0.000984627445307152
$ws.Cells.Item(4, 1).Value = 74
$ws.Cells.Item(4, 2).Value = "`$\langle qq \vert qq \rangle`$"
$ws.Cells.Item(4, 3).Value = 0.0007132819033761118
$ws.Cells.Item(5, 1).Value = 35
$ws.Cells.Item(5, 2).Value = "`$F_{q}^{\text{SCF}}`$"
$ws.Cells.Item(5, 3).Value = 0.0004343386760998269
$ws.Cells.Item(6, 1).Value = 46
$ws.Cells.Item(6, 2).Value = "`$\eta_{s}`$"
$ws.Cells.Item(6, 3).Value = 0.0004246014637988282
$ws.Cells.Item(7, 1).Value = 22
$ws.Cells.Item(7, 2).Value = "h`$_{s}`$"
$ws.Cells.Item(7, 3).Value = 0.0004201494299917982
$ws.Cells.Item(8, 1).Value = 63
$ws.Cells.Item(8, 2).Value = "`$(F_{p}^{\text{SCF}})_{3}`$"
$ws.Cells.Item(8, 3).Value = 0.0003118570237438018
$ws.Cells.Item(9, 1).Value = 19
$ws.Cells.Item(9, 2).Value = "h`$_{rs}^{1}`$"
$ws.Cells.Item(9, 3).Value = 0.0003035080493472498
$ws.Cells.Item(10, 1).Value = 13
$ws.Cells.Item(10, 2).Value = "h`$_{qs}`$"
$ws.Cells.Item(10, 3).Value = 0.0002862594811628542
$ws.Cells.Item(11, 1).Value = 62
$ws.Cells.Item(11, 2).Value = "`$(\eta_{r})_{2}`$"
$ws.Cells.Item(11, 3).Value = 0.0002581516543079358
$ws.Cells.Item(12, 1).Value = 55
$ws.Cells.Item(12, 2).Value = "`$(F_{p}^{\text{SCF}})_{2}`$"
$ws.Cells.Item(12, 3).Value = 0.0001759560042575764
$ws.Cells.Item(13, 1).Value = 16
$ws.Cells.Item(13, 2).Value = "h`$_{r}^{2}`$"
$ws.Cells.Item(13, 3).Value = 0.0001604164751221265
$ws.Cells.Item(14, 1).Value = 43
$ws.Cells.Item(14, 2).Value = "`$F_{s}^{\text{SCF}}`$"
$ws.Cells.Item(14, 3).Value = 0.0001546800836887528
$ws.Cells.Item(15, 1).Value = 67
$ws.Cells.Item(15, 2).Value = "`$(F_{r}^{\text{SCF}})_{3}`$"
$ws.Cells.Item(15, 3).Value = 0.0001375286233212047
$ws.Cells.Item(16, 1).Value = 3
$ws.Cells.Item(16, 2).Value = "h`$_{p}^{3}`$"
$ws.Cells.Item(16, 3).Value = 0.0001286380260166742
$ws.Cells.Item(17, 1).Value = 92
$ws.Cells.Item(17, 2).Value = "`$(\langle rr \vert rr \rangle)_{2}`$"
$ws.Cells.Item(17, 3).Value = 0.0001255336461939345
$ws.Cells.Item(18, 1).Value = 1
$ws.Cells.Item(18, 2).Value = "h`$_{p}^{1}`$"
$ws.Cells.Item(18, 3).Value = 0.0001162080009106264
$ws.Cells.Item(19, 1).Value = 96
$ws.Cells.Item(19, 2).Value = "`$(\langle rs \vert sr \rangle)_{2}`$"
$ws.Cells.Item(19, 3).Value = 0.0001147639773823396
$ws.Cells.Item(20, 1).Value = 23
$ws.Cells.Item(20, 2).Value = "typ_0"
$ws.Cells.Item(20, 3).Value = 0.0001105577958307923
$ws.Cells.Item(21, 1).Value = 11
$ws.Cells.Item(21, 2).Value = "h`$_{pr}^{3}`$"
$ws.Cells.Item(21, 3).Value = 0.000108790491467365
$ws.Cells.Item(22, 1).Value = 0
$ws.Cells.Item(22, 2).Value = "h`$_{p}^{0}`$"
$ws.Cells.Item(22, 3).Value = 0.0001028680679069994
$ws.Cells.Item(23, 1).Value = 18
$ws.Cells.Item(23, 2).Value = "h`$_{rs}^{0}`$"
$ws.Cells.Item(23, 3).Value = 0.00009949260314523085
$ws.Cells.Item(24, 1).Value = 93
$ws.Cells.Item(24, 2).Value = "`$(\langle pq \vert pq \rangle)_{2}`$"
$ws.Cells.Item(24, 3).Value = 0.00009922286855457385
$ws.Cells.Item(25, 1).Value = 20
$ws.Cells.Item(25, 2).Value = "h`$_{rs}^{2}`$"
$ws.Cells.Item(25, 3).Value = 0.00009861205476116337
$ws.Cells.Item(26, 1).Value = 76
$ws.Cells.Item(26, 2).Value = "`$\langle ss \vert ss \rangle`$"
$ws.Cells.Item(26, 3).Value = 0.00009590242302953046
$ws.Cells.Item(27, 1).Value = 104
$ws.Cells.Item(27, 2).Value = "`$(\langle rs \vert sr \rangle)_{3}`$"
$ws.Cells.Item(27, 3).Value = 0.00008731079635207966
$ws.Cells.Item(28, 1).Value = 69
$ws.Cells.Item(28, 2).Value = "`$(F_{r})_{3}`$"
$ws.Cells.Item(28, 3).Value = 0.0000842529590016625
$ws.Cells.Item(29, 1).Value = 29
$ws.Cells.Item(29, 2).Value = "FI`$_{qs}`$"
$ws.Cells.Item(29, 3).Value = 0.00008258255719271565
$ws.Cells.Item(30, 1).Value = 2
$ws.Cells.Item(30, 2).Value = "h`$_{p}^{2}`$"
$ws.Cells.Item(30, 3).Value = 0.00007732577947473127
$ws.Cells.Item(31, 1).Value = 99
$ws.Cells.Item(31, 2).Value = "`$(\langle pp \vert pp \rangle)_{3}`$"
$ws.Cells.Item(31, 3).Value = 0.00007680864729294616
$ws.Cells.Item(32, 1).Value = 91
$ws.Cells.Item(32, 2).Value = "`$(\langle pp \vert pp \rangle)_{2}`$"
$ws.Cells.Item(32, 3).Value = 0.00007613564257039784
$ws.Cells.Item(33, 1).Value = 17
$ws.Cells.Item(33, 2).Value = "h`$_{r}^{3}`$"
$ws.Cells.Item(33, 3).Value = 0.00006807622889993819
$ws.Cells.Item(34, 1).Value = 59
$ws.Cells.Item(34, 2).Value = "`$(F_{r}^{\text{SCF}})_{2}`$"
$ws.Cells.Item(34, 3).Value = 0.00006597161394624604
$ws.Cells.Item(35, 1).Value = 95
$ws.Cells.Item(35, 2).Value = "`$(\langle rs\vert rs \rangle)_{2}`$"
$ws.Cells.Item(35, 3).Value = 0.00006404086148789688
$ws.Cells.Item(36, 1).Value = 79
$ws.Cells.Item(36, 2).Value = "`$(\langle rs\vert rs \rangle)_{0}`$"
$ws.Cells.Item(36, 3).Value = 0.00006364928337697161
$ws.Cells.Item(37, 1).Value = 100
$ws.Cells.Item(37, 2).Value = "`$(\langle rr \vert rr \rangle)_{3}`$"
$ws.Cells.Item(37, 3).Value = 0.00006146369422016981
$ws.Cells.Item(38, 1).Value = 7
$ws.Cells.Item(38, 2).Value = "h`$_{pq}^{3}`$"
$ws.Cells.Item(38, 3).Value = 0.00005407766626921774
$ws.Cells.Item(39, 1).Value = 21
$ws.Cells.Item(39, 2).Value = "h`$_{rs}^{3}`$"
$ws.Cells.Item(39, 3).Value = 0.00005342958638665683
$ws.Cells.Item(40, 1).Value = 102
$ws.Cells.Item(40, 2).Value = "`$(\langle pq \vert qp \rangle)_{3}`$"
$ws.Cells.Item(40, 3).Value = 0.00005151417243060589
$ws.Cells.Item(41, 1).Value = 94
$ws.Cells.Item(41, 2).Value = "`$(\langle pq \vert qp \rangle)_{2}`$"
$ws.Cells.Item(41, 3).Value = 0.00004686833980635831
$ws.Cells.Item(42, 1).Value = 61
$ws.Cells.Item(42, 2).Value = "`$(F_{r})_{2}`$"
$ws.Cells.Item(42, 3).Value = 0.00004460199601233939
$ws.Cells.Item(43, 1).Value = 25
$ws.Cells.Item(43, 2).Value = "typ_2"
$ws.Cells.Item(43, 3).Value = 0.00003881498899773134
$ws.Cells.Item(44, 1).Value = 103
$ws.Cells.Item(44, 2).Value = "`$(\langle rs\vert rs \rangle)_{3}`$"
$ws.Cells.Item(44, 3).Value = 0.00003756550782864879
$ws.Cells.Item(45, 1).Value = 97
$ws.Cells.Item(45, 2).Value = "`$(\langle pq \vert rs \rangle)_{3}`$"
$ws.Cells.Item(45, 3).Value = 0.00003597129740228116
$ws.Cells.Item(46, 1).Value = 31
$ws.Cells.Item(46, 2).Value = "`$(F_{p}^{\text{SCF}})_{0}`$"
$ws.Cells.Item(46, 3).Value = 0.00003413299240274127
$ws.Cells.Item(47, 1).Value = 5
$ws.Cells.Item(47, 2).Value = "h`$_{pq}^{1}`$"
$ws.Cells.Item(47, 3).Value = 0.00003287857972125735
$ws.Cells.Item(48, 1).Value = 101
$ws.Cells.Item(48, 2).Value = "`$(\langle pq \vert pq \rangle)_{3}`$"
$ws.Cells.Item(48, 3).Value = 0.00003258901023940577
$ws.Cells.Item(49, 1).Value = 70
$ws.Cells.Item(49, 2).Value = "`$(\eta_{r})_{3}`$"
$ws.Cells.Item(49, 3).Value = 0.00003169248415369284
$ws.Cells.Item(50, 1).Value = 80
$ws.Cells.Item(50, 2).Value = "`$(\langle rs \vert sr \rangle)_{0}`$"
$ws.Cells.Item(50, 3).Value = 0.00003142121264564686
$ws.Cells.Item(51, 1).Value = 77
$ws.Cells.Item(51, 2).Value = "`$(\langle pq \vert pq \rangle)_{0}`$"
$ws.Cells.Item(51, 3).Value = 0.00003074880017948693
$ws.Cells.Item(52, 1).Value = 86
$ws.Cells.Item(52, 2).Value = "`$(\langle pq \vert qp \rangle)_{1}`$"
$ws.Cells.Item(52, 3).Value = 0.00003011555107273129
$ws.Cells.Item(53, 1).Value = 10
$ws.Cells.Item(53, 2).Value = "h`$_{pr}^{2}`$"
$ws.Cells.Item(53, 3).Value = 0.00002991116039952024
$ws.Cells.Item(54, 1).Value = 51
$ws.Cells.Item(54, 2).Value = "`$(F_{r}^{\text{SCF}})_{1}`$"
$ws.Cells.Item(54, 3).Value = 0.00002692062856001079
$ws.Cells.Item(55, 1).Value = 87
$ws.Cells.Item(55, 2).Value = "`$(\langle rs\vert rs \rangle)_{1}`$"
$ws.Cells.Item(55, 3).Value = 0.00002585664482959177
$ws.Cells.Item(56, 1).Value = 15
$ws.Cells.Item(56, 2).Value = "h`$_{r}^{1}`$"
$ws.Cells.Item(56, 3).Value = 0.0000253042796412504
$ws.Cells.Item(57, 1).Value = 8
$ws.Cells.Item(57, 2).Value = "h`$_{pr}^{0}`$"
$ws.Cells.Item(57, 3).Value = 0.00002443353238766979
$ws.Cells.Item(58, 1).Value = 84
$ws.Cells.Item(58, 2).Value = "`$(\langle rr \vert rr \rangle)_{1}`$"
$ws.Cells.Item(58, 3).Value = 0.00002417763751490489
$ws.Cells.Item(59, 1).Value = 6
$ws.Cells.Item(59, 2).Value = "h`$_{pq}^{2}`$"
$ws.Cells.Item(59, 3).Value = 0.00002260044660406281
$ws.Cells.Item(60, 1).Value = 85
$ws.Cells.Item(60, 2).Value = "`$(\langle pq \vert pq \rangle)_{1}`$"
$ws.Cells.Item(60, 3).Value = 0.00002137217709425837
$ws.Cells.Item(61, 1).Value = 45
$ws.Cells.Item(61, 2).Value = "`$F_{s}`$"
$ws.Cells.Item(61, 3).Value = 0.0000210607327728464
$ws.Cells.Item(62, 1).Value = 39
$ws.Cells.Item(62, 2).Value = "`$(F_{r}^{\text{SCF}})_{0}`$"
$ws.Cells.Item(62, 3).Value = 0.00001997739039343899
$ws.Cells.Item(63, 1).Value = 83
$ws.Cells.Item(63, 2).Value = "`$(\langle pp \vert pp \rangle)_{1}`$"
$ws.Cells.Item(63, 3).Value = 0.00001874288101563429
$ws.Cells.Item(64, 1).Value = 65
$ws.Cells.Item(64, 2).Value = "`$(F_{p})_{3}`$"
$ws.Cells.Item(64, 3).Value = 0.00001828192505932386
$ws.Cells.Item(65, 1).Value = 78
$ws.Cells.Item(65, 2).Value = "`$(\langle pq \vert qp \rangle)_{0}`$"
$ws.Cells.Item(65, 3).Value = 0.00001754658968561231
$ws.Cells.Item(66, 1).Value = 89
$ws.Cells.Item(66, 2).Value = "`$(\langle pq \vert rs \rangle)_{2}`$"
$ws.Cells.Item(66, 3).Value = 0.00001735712531551862
$ws.Cells.Item(67, 1).Value = 14
$ws.Cells.Item(67, 2).Value = "h`$_{r}^{0}`$"
$ws.Cells.Item(67, 3).Value = 0.00001579347350390622
$ws.Cells.Item(68, 1).Value = 4
$ws.Cells.Item(68, 2).Value = "h`$_{pq}^{0}`$"
$ws.Cells.Item(68, 3).Value = 0.00001562829254186308
$ws.Cells.Item(69, 1).Value = 26
$ws.Cells.Item(69, 2).Value = "typ_3"
$ws.Cells.Item(69, 3).Value = 0.00001465079096237254
$ws.Cells.Item(70, 1).Value = 88
$ws.Cells.Item(70, 2).Value = "`$(\langle rs \vert sr \rangle)_{1}`$"
$ws.Cells.Item(70, 3).Value = 0.00001422888378829175
$ws.Cells.Item(71, 1).Value = 75
$ws.Cells.Item(71, 2).Value = "`$(\langle rr \vert rr \rangle)_{0}`$"
$ws.Cells.Item(71, 3).Value = 0.00001146834853876853
$ws.Cells.Item(72, 1).Value = 30
$ws.Cells.Item(72, 2).Value = "FA`$_{qs}`$"
$ws.Cells.Item(72, 3).Value = 0.00001087359119852859
$ws.Cells.Item(73, 1).Value = 33
$ws.Cells.Item(73, 2).Value = "`$(F_{p})_{0}`$"
$ws.Cells.Item(73, 3).Value = 0.0000108044737063564
$ws.Cells.Item(74, 1).Value = 9
$ws.Cells.Item(74, 2).Value = "h`$_{pr}^{1}`$"
$ws.Cells.Item(74, 3).Value = 0.0000100406768882559
$ws.Cells.Item(75, 1).Value = 66
$ws.Cells.Item(75, 2).Value = "`$(\eta_{p})_{3}`$"
$ws.Cells.Item(75, 3).Value = 0.000009977273993727336
$ws.Cells.Item(76, 1).Value = 47
$ws.Cells.Item(76, 2).Value = "`$(F_{p}^{\text{SCF}})_{1}`$"
$ws.Cells.Item(76, 3).Value = 0.000009857977736798404
$ws.Cells.Item(77, 1).Value = 57
$ws.Cells.Item(77, 2).Value = "`$(F_{p})_{2}`$"
$ws.Cells.Item(77, 3).Value = 0.0000094347172948471
$ws.Cells.Item(78, 1).Value = 24
$ws.Cells.Item(78, 2).Value = "typ_1"
$ws.Cells.Item(78, 3).Value = 0.000008824820842308125
$ws.Cells.Item(79, 1).Value = 58
$ws.Cells.Item(79, 2).Value = "`$(\eta_{p})_{2}`$"
$ws.Cells.Item(79, 3).Value = 0.00000825424155796253
$ws.Cells.Item(80, 1).Value = 81
$ws.Cells.Item(80, 2).Value = "`$(\langle pq \vert rs \rangle)_{1}`$"
$ws.Cells.Item(80, 3).Value = 0.000007402657347032949
$ws.Cells.Item(81, 1).Value = 38
$ws.Cells.Item(81, 2).Value = "`$\eta_{q}`$"
$ws.Cells.Item(81, 3).Value = 0.000006476609795201717
$ws.Cells.Item(82, 1).Value = 41
$ws.Cells.Item(82, 2).Value = "`$(F_{r})_{0}`$"
$ws.Cells.Item(82, 3).Value = 0.000005653440425215473
$ws.Cells.Item(83, 1).Value = 53
$ws.Cells.Item(83, 2).Value = "`$(F_{r})_{1}`$"
$ws.Cells.Item(83, 3).Value = 0.000005525560050651195
$ws.Cells.Item(84, 1).Value = 73
$ws.Cells.Item(84, 2).Value = "`$(\langle pp \vert pp \rangle)_{0}`$"
$ws.Cells.Item(84, 3).Value = 0.000002933309356432915
$ws.Cells.Item(85, 1).Value = 71
$ws.Cells.Item(85, 2).Value = "`$(\langle pq \vert rs \rangle)_{0}`$"
$ws.Cells.Item(85, 3).Value = 0.00000254789375209164
$ws.Cells.Item(86, 1).Value = 50
$ws.Cells.Item(86, 2).Value = "`$(\eta_{p})_{1}`$"
$ws.Cells.Item(86, 3).Value = 0.000002359844099875356
$ws.Cells.Item(87, 1).Value = 34
$ws.Cells.Item(87, 2).Value = "`$(\eta_{p})_{0}`$"
$ws.Cells.Item(87, 3).Value = 0.000002348289435587717
$ws.Cells.Item(88, 1).Value = 49
$ws.Cells.Item(88, 2).Value = "`$(F_{p})_{1}`$"
$ws.Cells.Item(88, 3).Value = 0.000001387510388167063
$ws.Cells.Item(89, 1).Value = 42
$ws.Cells.Item(89, 2).Value = "`$(\eta_{r})_{0}`$"
$ws.Cells.Item(89, 3).Value = 0.0000005683718546081586
$ws.Cells.Item(90, 1).Value = 54
$ws.Cells.Item(90, 2).Value = "`$(\eta_{r})_{1}`$"
$ws.Cells.Item(90, 3).Value = 0.0000005315252818727352
$ws.Cells.Item(91, 1).Value = 82
$ws.Cells.Item(91, 2).Value = "`$(\langle pq \vert sr \rangle)_{1}`$"
$ws.Cells.Item(91, 3).Value = 0.00000007105690486168858
$ws.Cells.Item(92, 1).Value = 98
$ws.Cells.Item(92, 2).Value = "`$(\langle pq \vert sr \rangle)_{3}`$"
$ws.Cells.Item(92, 3).Value = 0.00000003179084676730028
$ws.Cells.Item(93, 1).Value = 90
$ws.Cells.Item(93, 2).Value = "`$(\langle pq \vert sr \rangle)_{2}`$"
$ws.Cells.Item(93, 3).Value = 0.00000003151891295082749
$ws.Cells.Item(94, 1).Value = 72
$ws.Cells.Item(94, 2).Value = "`$(\langle pq \vert sr \rangle)_{0}`$"
$ws.Cells.Item(94, 3).Value = 0.0000000293781430289728
$ws.Cells.Item(95, 1).Value = 44
$ws.Cells.Item(95, 2).Value = "`$\omega_{s}`$"
$ws.Cells.Item(95, 3).Value = 0.00000002212451673699822
$ws.Cells.Item(96, 1).Value = 27
$ws.Cells.Item(96, 2).Value = "`$\mathbf{b}`$"
$ws.Cells.Item(96, 3).Value = 0.00000001458872846755756
$ws.Cells.Item(97, 1).Value = 28
$ws.Cells.Item(97, 2).Value = "F`$_{qs}`$"
$ws.Cells.Item(97, 3).Value = 0.00000001328955405927034
$ws.Cells.Item(98, 1).Value = 40
$ws.Cells.Item(98, 2).Value = "`$(\omega_{r})_{0}`$"
$ws.Cells.Item(98, 3).Value = 0.00000001195765122952185
$ws.Cells.Item(99, 1).Value = 32
$ws.Cells.Item(99, 2).Value = "`$(\omega_{p})_{0}`$"
$ws.Cells.Item(99, 3).Value = 0.000000008292245709624312
$ws.Cells.Item(100, 1).Value = 48
$ws.Cells.Item(100, 2).Value = "`$(\omega_{p})_{1}`$"
$ws.Cells.Item(100, 3).Value = 0.000000008189232610187569
$ws.Cells.Item(101, 1).Value = 64
$ws.Cells.Item(101, 2).Value = "`$(\omega_{p})_{3}`$"
$ws.Cells.Item(101, 3).Value = 0.000000006968103791420767
$ws.Cells.Item(102, 1).Value = 68
$ws.Cells.Item(102, 2).Value = "`$(\omega_{r})_{3}`$"
$ws.Cells.Item(102, 3).Value = 0.000000006959363277441956
$ws.Cells.Item(103, 1).Value = 56
$ws.Cells.Item(103, 2).Value = "`$(\omega_{p})_{2}`$"
$ws.Cells.Item(103, 3).Value = 0.00000000629978959740921
$ws.Cells.Item(104, 1).Value = 60
$ws.Cells.Item(104, 2).Value = "`$(\omega_{r})_{2}`$"
$ws.Cells.Item(104, 3).Value = 0.000000006286472104813377
$ws.Cells.Item(105, 1).Value = 52
$ws.Cells.Item(105, 2).Value = "`$(\omega_{r})_{1}`$"
$ws.Cells.Item(105, 3).Value = 0.000000004509072687091276
$ws.Cells.Item(106, 1).Value = 36
$ws.Cells.Item(106, 2).Value = "`$\omega_{q}`$"
$ws.Cells.Item(106, 3).Value = 0.000000001118169166147726
